$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.213243961334229
$ws.Range("B1").Value = 2.443443059921265
$ws.Range("C1").Value = 4.771324634552002
$ws.Range("D1").Value = 2.523614883422852
$ws.Range("E1").Value = 1.081949591636658
